# Penalty Reward System (unfinished) - update forecast week dates and
# "actual"/MyForecast values by one week, and refresh dependent Summary stats.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": shift Week_Start_Date (col B) forward one
# week and update MyForecast (col D) values for rows 2-17. Columns B keep
# plain text (YYYY-MM-DD) formatting just like the source data, so force the
# cell format to Text before assigning the date strings to avoid Excel
# auto-converting them into date serials.

$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "2025-01-12"
$ws1.Range("D2").Value = 26

$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "2025-01-19"
$ws1.Range("D3").Value = 29

$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "2025-01-26"
$ws1.Range("D4").Value = 33

$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "2025-02-02"
$ws1.Range("D5").Value = 37

$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = "2025-02-09"
$ws1.Range("D6").Value = 40

$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2025-02-16"
$ws1.Range("D7").Value = 43

$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "2025-02-23"
$ws1.Range("D8").Value = 31

$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "2025-03-02"
$ws1.Range("D9").Value = 31

$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "2025-03-09"
$ws1.Range("D10").Value = 28

$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = "2025-03-16"
$ws1.Range("D11").Value = 28

$ws1.Range("B12").NumberFormat = "@"
$ws1.Range("B12").Value = "2025-03-23"
$ws1.Range("D12").Value = 27

$ws1.Range("B13").NumberFormat = "@"
$ws1.Range("B13").Value = "2025-03-30"
$ws1.Range("D13").Value = 27

$ws1.Range("B14").NumberFormat = "@"
$ws1.Range("B14").Value = "2025-04-06"
$ws1.Range("D14").Value = 24

$ws1.Range("B15").NumberFormat = "@"
$ws1.Range("B15").Value = "2025-04-13"
$ws1.Range("D15").Value = 31

$ws1.Range("B16").NumberFormat = "@"
$ws1.Range("B16").Value = "2025-04-20"
$ws1.Range("D16").Value = 31

$ws1.Range("B17").NumberFormat = "@"
$ws1.Range("B17").Value = "2025-04-27"
$ws1.Range("D17").Value = 22

# --- Sheet "Summary": refresh the metrics that depend on the shifted
# historical/forecast window. All values on this sheet are stored as plain
# text, so force Text format before writing the (numeric-looking) strings.

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "2023-01-22 to 2025-01-05"

$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "16"

$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = "1364 units"

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "488"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "269"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "125"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "43"

$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "2025-04-27"
